# Data file and Registration test case added
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the registration test data rows (values only - header row text is unchanged)
$ws.Range("A2").Value = "FirstName51"
$ws.Range("B2").Value = "LastName51"
$ws.Range("C2").Value = "ncitester51@nih.gov"

$ws.Range("A3").Value = "FirstName52"
$ws.Range("B3").Value = "LastName52"
$ws.Range("C3").Value = "ncitester52@nih.gov"

$ws.Range("A4").Value = "FirstName53"
$ws.Range("B4").Value = "LastName53"
$ws.Range("C4").Value = "ncitester53@nih.gov"

$ws.Range("A5").Value = "FirstName54"
$ws.Range("B5").Value = "LastName54"
$ws.Range("C5").Value = "ncitester54@nih.gov"

# Reset row heights back to default (no more custom/tall rows) - AutoFit
# clears the custom-height flag the same way Excel does when a row reverts
# to the sheet's standard height.
$ws.Rows.Item(1).AutoFit()
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(5).AutoFit()

# Set explicit column widths for A:E (values chosen so the saved character
# width matches the target file as closely as this engine's width grid allows)
$ws.Columns.Item(1).ColumnWidth = 14.833333333333334
$ws.Columns.Item(2).ColumnWidth = 17.166666666666668
$ws.Columns.Item(3).ColumnWidth = 35.5
$ws.Columns.Item(4).ColumnWidth = 32.666666666666664
$ws.Columns.Item(5).ColumnWidth = 19.666666666666668

# Move the selection to A2
$ws.Range("A2").Select()
